$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "Update index.py"
$ws.Range("C3").Value = "riya-morankar"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "from edit1 to main"

# The target date ("2025-06-17") must be stored as literal text, not an
# auto-converted date serial number. Prefixing with an apostrophe is the
# standard Excel mechanism to force text entry (mirrors what a user would
# type into the Ribbon/formula bar to keep a date-shaped string as text).
$ws.Range("F3").Value = "'2025-06-17"
